$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.211.36'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '2.899.55'
$ws.Range('E3').Value = '  +3.85%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'352.69"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').Value = "'113.09"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.75%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').Value = "'39.99"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = "'0.0861"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('D13').Value = "'19.86"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = "'7.74"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('D15').Value = '3.358.90'
$ws.Range('E15').Value = '  +3.93%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = "'0.994"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.55%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.881.41'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').Value = '52.218.71'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').Value = "'3.30"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.46%  '
$ws.Range('D21').Value = "'14.17"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.95%  '
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Value = "'70.85"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.02%  '
$ws.Range('D24').Value = "'268.85"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E25').Value = '  +1.52%  '
$ws.Range('E26').Value = '  +8.34%  '
$ws.Range('E27').Value = '  +2.78%  '
$ws.Range('D28').Value = "'1.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('D30').Value = "'0.102"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +15.47%  '
$ws.Range('D31').Value = "'6.61"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.34%  '
$ws.Range('D32').Value = "'37.51"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.46%  '
$ws.Range('D33').Value = "'6.20"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.13%  '
$ws.Range('D34').Value = "'53.17"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.11%  '
$ws.Range('D35').Value = "'0.0448"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range('E36').Value = '  -13.14%  '
$ws.Range('D37').Value = "'0.998"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = "'3.32"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.54%  '
$ws.Range('D39').Value = "'18.84"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.60%  '
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('E41').Value = '  +10.75%  '
$ws.Range('E42').Value = '  +1.85%  '
$ws.Range('D43').Value = "'23.01"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.55%  '
$ws.Range('D44').Value = "'2.60"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.61%  '
$ws.Range('D45').Value = "'119.85"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('E46').Value = '  -1.49%  '
$ws.Range('E47').Value = '  +3.18%  '
$ws.Range('D48').Value = '2.174.98'
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('D49').Value = "'0.262"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +21.58%  '
$ws.Range('D50').Value = "'0.0348"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.16%  '
$ws.Range('E51').Value = '  -0.10%  '
